$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns present in the data rows (B-F, G, I-N). Column G stays at 1 (unchanged)
# and column H is unused, so we only write to B,C,D,E,F,I,J,K,L,M,N.
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New values for rows 2..25 (one array per row, matching $cols order)
$data = @(
    @(1.02, 1.025522496608464, 1.03066626895736, 1.048619401804276, 1.052975150287387, 1.032632857840732, 1.030691046903377, 1.033476619816118, 1.051378741284249, 1.055722397776916, 1.014279502194414),
    @(1.02, 1.026274657509716, 1.031219440644416, 1.049804016735474, 1.054209044575133, 1.032765141769633, 1.031083595773807, 1.033838879628324, 1.052374474540896, 1.056768165742035, 1.014409663098113),
    @(1.02, 1.026761694012317, 1.031577615346163, 1.050571574866488, 1.055008516555137, 1.032849690326169, 1.031337256182172, 1.034072815488076, 1.053019239593225, 1.05744534237714, 1.014493757343662),
    @(1.02, 1.026966523605136, 1.031728246814353, 1.050894503980588, 1.05534486862741, 1.032884983082129, 1.031443811488675, 1.034171048457428, 1.053290408984004, 1.057730146073181, 1.014529079446749),
    @(1.02, 1.027000920010404, 1.031753541669332, 1.050948739734295, 1.055401358527529, 1.03289089412914, 1.031461697667822, 1.034187535509634, 1.053335945985794, 1.057777972829232, 1.01453500834864),
    @(1.02, 1.026764430645062, 1.031579627877796, 1.05057588888745, 1.055013009909323, 1.032850162898461, 1.031338680308791, 1.034074128529753, 1.053022862538687, 1.057449147471664, 1.014494229442031),
    @(1.02, 1.025776622000803, 1.03085316617721, 1.049019534279231, 1.053391931374857, 1.032677780466816, 1.030823781490302, 1.033599144005392, 1.051715159107542, 1.056075717601364, 1.01432351709702),
    @(1.02, 1.024038632487922, 1.029574922483181, 1.046284934469218, 1.050543486943659, 1.03236601926295, 1.029913863701649, 1.032758603136234, 1.049414322491522, 1.053659340624081, 1.01402173046928),
    @(1.02, 1.022881846831755, 1.028724114678334, 1.044467151794706, 1.048649949625134, 1.032152832428658, 1.029305559443032, 1.032195912154671, 1.047882768524748, 1.052050952021864, 1.013819907515253),
    @(1.02, 1.022381408898359, 1.028356046666514, 1.043681278815207, 1.04783130816894, 1.032059259201236, 1.029041767242806, 1.031951720883916, 1.047220138234264, 1.051355098208609, 1.013732370275914),
    @(1.02, 1.022195594236362, 1.028219382107417, 1.043389556104418, 1.047527418899382, 1.032024312895603, 1.028943725036807, 1.031860936912173, 1.046974088870725, 1.051096715206761, 1.013699833359963),
    @(1.02, 1.022235448908636, 1.028248694727437, 1.043452123231675, 1.047592595562851, 1.03203181754258, 1.028964758034325, 1.031880414006922, 1.047026863635645, 1.051152135258177, 1.013706813614964),
    @(1.02, 1.022366047958444, 1.028344748855513, 1.04365716114429, 1.047806184701973, 1.032056374381957, 1.02903366422209, 1.031944218287967, 1.047199798067007, 1.05133373839448, 1.013729681204792),
    @(1.02, 1.022446523698737, 1.028403937950566, 1.043783516331905, 1.047937809261026, 1.032071479627564, 1.029076111935608, 1.031983519580001, 1.047306359400101, 1.05144564174062, 1.013743767817722),
    @(1.02, 1.022915069022417, 1.028748549385093, 1.044519333644841, 1.048704306926075, 1.032159016035009, 1.029323058273046, 1.032212106985465, 1.04792675650549, 1.052097145923195, 1.013825714013906),
    @(1.02, 1.02320909875702, 1.028964806491943, 1.044981223652014, 1.049185450978422, 1.032213587844315, 1.029477856769984, 1.032355349342369, 1.048316060154224, 1.052505974700067, 1.013877077692258),
    @(1.02, 1.02338064540736, 1.029090978048991, 1.045250755857558, 1.049466216999948, 1.032245296863581, 1.029568110187494, 1.032438847787075, 1.048543186780407, 1.05274449447231, 1.013907023092843),
    @(1.02, 1.023439145813917, 1.029134004743904, 1.045342679676198, 1.049561971810889, 1.032256088144745, 1.029598877809388, 1.032467309680074, 1.048620639992842, 1.052825833208466, 1.013917231287996),
    @(1.02, 1.023177547579022, 1.028941600783469, 1.044931654861325, 1.049133816073728, 1.032207745399396, 1.029461252272531, 1.032339986205546, 1.048274286141165, 1.052462105362632, 1.013871568318314),
    @(1.02, 1.022327587841767, 1.028316461852625, 1.043596777473623, 1.047743282766806, 1.03204914821893, 1.02901337466195, 1.031925431730398, 1.047148870954418, 1.051280258368816, 1.013722947865565),
    @(1.02, 1.02179359170223, 1.027923715510285, 1.042758560047488, 1.046870101853027, 1.031948338405777, 1.02873144090084, 1.031664319859697, 1.046441746330358, 1.050537693703513, 1.013629378974199),
    @(1.02, 1.02207663406097, 1.028131888493373, 1.04320281346728, 1.047332887261256, 1.032001883060763, 1.028880930799735, 1.031802783940124, 1.046816562223371, 1.050931293128649, 1.013678993381338),
    @(1.02, 1.023191804062218, 1.028952086344197, 1.044954052494445, 1.049157147284768, 1.032210385725152, 1.029468755247484, 1.032346928307796, 1.048293161857418, 1.052481927851979, 1.013874057811052),
    @(1.02, 1.024487620613445, 1.029905147233473, 1.046990960958669, 1.051278920189581, 1.032447561745661, 1.030149401645568, 1.032976318339297, 1.050008730631638, 1.054283585219596, 1.01409986238528)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "$($cols[$j])$row"
        $ws.Range($addr).Value = $vals[$j]
    }
}
